$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.095.87"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.64"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.48"
$ws.Range("E5").Value = "  -3.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4637"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07866"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9610"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.93"
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.689"
$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.914"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.763.08"
$ws.Range("E14").Value = "  -7.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06846"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.27"
$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009926"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.65"
$ws.Range("E19").Value = "  -2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.108.77"
$ws.Range("E21").Value = "  -2.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.330"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  -2.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.093"
$ws.Range("E24").Value = "  -1.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062.97"
$ws.Range("E25").Value = "  -3.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.19"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.13"
$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.701"
$ws.Range("E28").Value = "  -7.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.967"
$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.06"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9370"
$ws.Range("E31").Value = "  -4.03%  "

$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.272"
$ws.Range("E33").Value = "  -1.81%  "

$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.296"
$ws.Range("E35").Value = "  -5.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05851"
$ws.Range("E36").Value = "  -5.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02128"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.144"
$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.776"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5598"
$ws.Range("E40").Value = "  -2.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.902"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("E42").Value = "  -1.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07240"
$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.62"
$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5262"
$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.129"
$ws.Range("E46").Value = "  -10.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.117"
$ws.Range("E47").Value = "  -10.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.836"
$ws.Range("E48").Value = "  -4.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.71"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.324"
$ws.Range("E51").Value = "  +0.19%  "
